$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A1:A46")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1")) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

$sortRange.RemoveDuplicates(1)

$ws.Range("A27").Select()
